# Update the table style applied to every table in the deck: swap the
# deck's custom default table style ("Table_0",
# {5408B22E-FE57-4403-B3ED-57E334498C85}) for the built-in table style
# {D8F913B0-9F57-40F3-B3C3-D13D074FBAF3} (as picked from the Table
# Styles gallery on the Table Design ribbon).
#
# The presentation has three tables (one each on three slides), all
# currently carrying the old style id. Rather than hard-code slide
# numbers, walk every slide/shape and re-style whichever shapes
# actually contain a table using that old style.

$oldStyleId = "{5408B22E-FE57-4403-B3ED-57E334498C85}"
$newStyleId = "{D8F913B0-9F57-40F3-B3C3-D13D074FBAF3}"

$p = $ppt.ActivePresentation

for ($slideIndex = 1; $slideIndex -le $p.Slides.Count; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)

    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
